# Update the multiplication problems shown in the worksheet table.
# Each cell contains a single literal "AxB=" expression; replace the
# old expression text with the new one. Every "old" string is unique
# within the document, so a simple Find/Replace (wrap = none,
# MatchWholeWord not needed since '=' delimits the match) is safe.
$d = $word.ActiveDocument

$d.Content.Find.Execute("96×26=", $true, $false, $false, $false, $false, $true, 1, $false, "77×26=", 2) | Out-Null
$d.Content.Find.Execute("66×39=", $true, $false, $false, $false, $false, $true, 1, $false, "20×56=", 2) | Out-Null
$d.Content.Find.Execute("31×16=", $true, $false, $false, $false, $false, $true, 1, $false, "87×26=", 2) | Out-Null
$d.Content.Find.Execute("81×94=", $true, $false, $false, $false, $false, $true, 1, $false, "27×38=", 2) | Out-Null
$d.Content.Find.Execute("77×63=", $true, $false, $false, $false, $false, $true, 1, $false, "97×53=", 2) | Out-Null
$d.Content.Find.Execute("58×73=", $true, $false, $false, $false, $false, $true, 1, $false, "77×81=", 2) | Out-Null
$d.Content.Find.Execute("63×65=", $true, $false, $false, $false, $false, $true, 1, $false, "28×66=", 2) | Out-Null
$d.Content.Find.Execute("41×95=", $true, $false, $false, $false, $false, $true, 1, $false, "51×45=", 2) | Out-Null
$d.Content.Find.Execute("65×32=", $true, $false, $false, $false, $false, $true, 1, $false, "24×12=", 2) | Out-Null
$d.Content.Find.Execute("31×81=", $true, $false, $false, $false, $false, $true, 1, $false, "24×94=", 2) | Out-Null
$d.Content.Find.Execute("30×80=", $true, $false, $false, $false, $false, $true, 1, $false, "39×57=", 2) | Out-Null
$d.Content.Find.Execute("15×91=", $true, $false, $false, $false, $false, $true, 1, $false, "20×52=", 2) | Out-Null
$d.Content.Find.Execute("76×45=", $true, $false, $false, $false, $false, $true, 1, $false, "97×45=", 2) | Out-Null
$d.Content.Find.Execute("64×42=", $true, $false, $false, $false, $false, $true, 1, $false, "59×60=", 2) | Out-Null
$d.Content.Find.Execute("61×51=", $true, $false, $false, $false, $false, $true, 1, $false, "78×17=", 2) | Out-Null
$d.Content.Find.Execute("22×59=", $true, $false, $false, $false, $false, $true, 1, $false, "15×68=", 2) | Out-Null
$d.Content.Find.Execute("42×92=", $true, $false, $false, $false, $false, $true, 1, $false, "20×75=", 2) | Out-Null
$d.Content.Find.Execute("84×55=", $true, $false, $false, $false, $false, $true, 1, $false, "15×45=", 2) | Out-Null
$d.Content.Find.Execute("47×17=", $true, $false, $false, $false, $false, $true, 1, $false, "25×35=", 2) | Out-Null
$d.Content.Find.Execute("50×92=", $true, $false, $false, $false, $false, $true, 1, $false, "79×32=", 2) | Out-Null
$d.Content.Find.Execute("29×82=", $true, $false, $false, $false, $false, $true, 1, $false, "36×41=", 2) | Out-Null
$d.Content.Find.Execute("53×63=", $true, $false, $false, $false, $false, $true, 1, $false, "22×55=", 2) | Out-Null
$d.Content.Find.Execute("50×23=", $true, $false, $false, $false, $false, $true, 1, $false, "71×67=", 2) | Out-Null
$d.Content.Find.Execute("37×97=", $true, $false, $false, $false, $false, $true, 1, $false, "31×47=", 2) | Out-Null
$d.Content.Find.Execute("17×88=", $true, $false, $false, $false, $false, $true, 1, $false, "16×21=", 2) | Out-Null
